# RDCC-5182 Added Version check
# Adds a new "VERSION" worksheet (placed after the existing sheet) containing
# a small "File Version" / "vx.xx" note in A6/B6, and makes it the active tab.

$wb = $excel.ActiveWorkbook
$mainSheet = $wb.Worksheets.Item(1)

# Add the new sheet right after the existing "Service to CW Roles Mapping" sheet.
$versionSheet = $wb.Worksheets.Add($null, $mainSheet)
$versionSheet.Name = "VERSION"

# Populate the version info starting at row 6 (matches the target layout).
$versionSheet.Range("A6").Value = "File Version"
$versionSheet.Range("B6").Value = "vx.xx"

# Make B6 the active selection on the new sheet, and make the new sheet the
# active/selected tab (mirrors tabSelected moving from sheet1 to sheet2).
$versionSheet.Range("B6").Select()
